$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B10 was previously stored as a text value "20" (inline string); it
# should become a real numeric 20.
$ws.Range("B10").Value = 20

# Append new attendance row 11: Sina Neak / 21 / Male / Kompot / Class B 2025 / image
$ws.Range("A11").Value = "Sina Neak"

# B11's "21" must stay a text value (like the other Age cells originally
# were) rather than being auto-coerced into a number by a plain .Value
# assignment. Write it through a Text-formatted scratch cell and paste
# only the values across, then clean the scratch cell back up so it
# doesn't leave stray formatting or widen the used range.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "21"
$ws.Range("Z1").Copy()
$ws.Range("B11").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("C11").Value = "Male"
$ws.Range("D11").Value = "Kompot"
$ws.Range("E11").Value = "Class B 2025"
$ws.Range("F11").Value = "image\34e67118c6594bcd88ac6902475643c7.png"
